$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# A new weekly record is inserted at row 3, pushing every existing record
# (previously rows 3-68) down by one row (new rows 4-69).
$ws.Rows.Item(3).Insert()

# Populate the newly inserted row 3 with the new observation. All the
# "static" columns (mercado/region/categoria/etc.) are identical across
# every record in this sheet, so reuse those values; only the date and the
# price/volume figures (D, J, K, L, M, P) are new for this record.
$ws.Range("A3").Value = 2
$ws.Range("B3").Value = "Comercializadora del Agro de Limarí"
$ws.Range("C3").Value = "Coquimbo"
$ws.Range("D3").Value = 44812
$ws.Range("E3").Value = 4
$ws.Range("F3").Value = 100112026
$ws.Range("G3").Value = "Haba"
$ws.Range("H3").Value = "Sin especificar"
$ws.Range("I3").Value = "Primera"
$ws.Range("J3").Value = 600
$ws.Range("K3").Value = 5000
$ws.Range("L3").Value = 6000
$ws.Range("M3").Value = 5500
$ws.Range("N3").Value = "$/saco 25 kilos"
$ws.Range("O3").Value = "Provincia de Limarí"
$ws.Range("P3").Value = 220
$ws.Range("Q3").Value = 25
$ws.Range("R3").Value = "Hortaliza"
